$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-08-02 Friday" "2024-08-03 Saturday"

Replace-Text "510÷7=" "104÷4="
Replace-Text "134÷8=" "518÷4="
Replace-Text "930÷5=" "351÷3="
Replace-Text "205÷3=" "151÷6="
Replace-Text "109÷4=" "843÷2="

Replace-Text "753÷6=" "614÷3="
Replace-Text "848÷3=" "351÷3="
Replace-Text "340÷7=" "719÷8="
Replace-Text "369÷6=" "746÷4="
Replace-Text "950÷8=" "610÷4="

Replace-Text "793÷7=" "325÷7="
Replace-Text "538÷3=" "827÷2="
Replace-Text "928÷6=" "383÷4="
Replace-Text "663÷6=" "559÷3="
Replace-Text "415÷2=" "226÷2="

Replace-Text "712÷4=" "565÷5="
Replace-Text "744÷8=" "265÷9="
Replace-Text "536÷9=" "834÷2="
Replace-Text "595÷2=" "365÷5="
Replace-Text "141÷6=" "642÷4="

Replace-Text "711÷7=" "810÷3="
Replace-Text "879÷7=" "890÷5="
Replace-Text "855÷2=" "750÷4="
Replace-Text "245÷7=" "584÷4="
Replace-Text "742÷4=" "987÷8="
